$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0021749262078608
$ws.Range("C2").Value = 0.0125835016311947
$ws.Range("D2").Value = 0.888768059655119
$ws.Range("E2").Value = 0.0214385583346279
$ws.Range("F2").Value = 0.0010874631039304
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.000466055615970172
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0.00869970483144322
$ws.Range("K2").Value = 0.997514370048159
$ws.Range("L2").Value = 0.00062140748796023
$ws.Range("M2").Value = 0.99906788876806
$ws.Range("N2").Value = 0.997980425664129
$ws.Range("O2").Value = 0.956190772098804
$ws.Range("P2").Value = 0.000310703743980115
$ws.Range("Q2").Value = 0.000310703743980115
$ws.Range("R2").Value = 0.000310703743980115
$ws.Range("S2").Value = 0.000310703743980115
$ws.Range("T2").Value = 0.997825073792139
$ws.Range("U2").Value = 0.105949976697219
$ws.Range("V2").Value = 0.024390243902439
$ws.Range("W2").Value = 0.0365076899176635
$ws.Range("X2").Value = 0.977784682305422

$ws.Range("B3").Value = 0.995960851328259
$ws.Range("C3").Value = 0.959453161410595
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.000466055615970172
$ws.Range("F3").Value = 0.92014913779711
$ws.Range("G3").Value = 0.000932111231940345
$ws.Range("H3").Value = 0.000310703743980115
$ws.Range("I3").Value = 0.0083890010874631
$ws.Range("J3").Value = 0.000155351871990057
$ws.Range("K3").Value = 0.000310703743980115
$ws.Range("L3").Value = 0.0428771166692559
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.000155351871990057
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0.00062140748796023
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.0772098803790586
$ws.Range("S3").Value = 0.99953394438403
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.858319092745068
$ws.Range("W3").Value = 0.000776759359950287
$ws.Range("X3").Value = 0.00201957433587075

$ws.Range("B4").Value = 0.00062140748796023
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.000155351871990057
$ws.Range("E4").Value = 0.974677644865621
$ws.Range("F4").Value = 0.00155351871990057
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0.99922324064005
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0.988348609600746
$ws.Range("K4").Value = 0.00201957433587075
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0.000466055615970172
$ws.Range("N4").Value = 0.0010874631039304
$ws.Range("O4").Value = 0.0400807829734348
$ws.Range("P4").Value = 0.000310703743980115
$ws.Range("Q4").Value = 0.99968929625602
$ws.Range("R4").Value = 0.00062140748796023
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0.00155351871990057
$ws.Range("U4").Value = 0.887214540935218
$ws.Range("V4").Value = 0.0414789498213453
$ws.Range("W4").Value = 0.00310703743980115
$ws.Range("X4").Value = 0.0167780021749262

$ws.Range("B5").Value = 0.000932111231940345
$ws.Range("C5").Value = 0.0209725027186578
$ws.Range("D5").Value = 0.10796955103309
$ws.Range("E5").Value = 0.00201957433587075
$ws.Range("F5").Value = 0.0705297498834861
$ws.Range("G5").Value = 0.99860183315209
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.989280720832686
$ws.Range("J5").Value = 0.000310703743980115
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0.94547149293149
$ws.Range("M5").Value = 0.000310703743980115
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.000310703743980115
$ws.Range("P5").Value = 0.99875718502408
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.915022526021439
$ws.Range("S5").Value = 0.000155351871990057
$ws.Range("T5").Value = 0.000466055615970172
$ws.Range("U5").Value = 0.000310703743980115
$ws.Range("V5").Value = 0.064781730619854
$ws.Range("W5").Value = 0.954637253378903
$ws.Range("X5").Value = 0.000466055615970172

